$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sending-cluster added ("ECs"). The LR-pair table now covers every
# combination of 3 sending clusters (ECs, FAPs, MuSCs) x 3 target clusters
# (FAPs, MuSCs, Resolving-Mac) for the Efnb3 -> Ephb2 pair, recomputed with
# the refreshed TPM values. Rewrite the data block (rows 2-10) in place.

$data = @(
    @('ECs', 'Efnb3', 'Ephb2', 'FAPs', 2, 0.6666666666666666, 0.2151756666666667, 0.645527, 0.1791915537270423, 0.1791915537270423, 3, 1, 8.882730333333333, 26.648191, 0.9347132976570145, 0.9347132976570145, 1.911347421295222, 17.202126791657, 0.1674927280964878, 0.1674927280964878),
    @('ECs', 'Efnb3', 'Ephb2', 'MuSCs', 2, 0.6666666666666666, 0.2151756666666667, 0.645527, 0.1791915537270423, 0.1791915537270423, 3, 1, 0.616144, 1.848432, 0.06483569448352988, 0.0648356944835299, 0.1325791959626667, 1.193212763664, 0.01161800883147554, 0.01161800883147555),
    @('ECs', 'Efnb3', 'Ephb2', 'Resolving-Mac', 2, 0.6666666666666666, 0.2151756666666667, 0.645527, 0.1791915537270423, 0.1791915537270423, 2, 0.6666666666666666, 0.004286, 0.012858, 0.0004510078594555965, 0.0004510078594555965, 0.0009222429073333332, 0.008300186166, 0.00008081679907895586, 0.00008081679907895587),
    @('FAPs', 'Efnb3', 'Ephb2', 'FAPs', 2, 0.6666666666666666, 0.5010686666666667, 1.503206, 0.4172742870736815, 0.4172742870736815, 3, 1, 8.882730333333333, 26.648191, 0.9347132976570145, 0.9347132976570145, 4.450857844482889, 40.05772060034599, 0.3900318248981205, 0.3900318248981205),
    @('FAPs', 'Efnb3', 'Ephb2', 'MuSCs', 2, 0.6666666666666666, 0.5010686666666667, 1.503206, 0.4172742870736815, 0.4172742870736815, 3, 1, 0.616144, 1.848432, 0.06483569448352988, 0.0648356944835299, 0.3087304525546667, 2.778574072992, 0.02705426819254195, 0.02705426819254196),
    @('FAPs', 'Efnb3', 'Ephb2', 'Resolving-Mac', 2, 0.6666666666666666, 0.5010686666666667, 1.503206, 0.4172742870736815, 0.4172742870736815, 2, 0.6666666666666666, 0.004286, 0.012858, 0.0004510078594555965, 0.0004510078594555965, 0.002147580305333333, 0.019328222748, 0.0001881939830189611, 0.0001881939830189612),
    @('MuSCs', 'Efnb3', 'Ephb2', 'FAPs', 2, 0.6666666666666666, 0.4845693333333334, 1.453708, 0.4035341591992763, 0.4035341591992763, 3, 1, 8.882730333333333, 26.648191, 0.9347132976570145, 0.9347132976570145, 4.304298715803111, 38.738688442228, 0.3771887446624062, 0.3771887446624062),
    @('MuSCs', 'Efnb3', 'Ephb2', 'MuSCs', 2, 0.6666666666666666, 0.4845693333333334, 1.453708, 0.4035341591992763, 0.4035341591992763, 3, 1, 0.616144, 1.848432, 0.06483569448352988, 0.0648356944835299, 0.2985644873173334, 2.687080385856, 0.02616341745951238, 0.02616341745951239),
    @('MuSCs', 'Efnb3', 'Ephb2', 'Resolving-Mac', 2, 0.6666666666666666, 0.4845693333333334, 1.453708, 0.4035341591992763, 0.4035341591992763, 2, 0.6666666666666666, 0.004286, 0.012858, 0.0004510078594555965, 0.0004510078594555965, 0.002076864162666666, 0.018691777464, 0.0001819970773576795, 0.0001819970773576795)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}
